$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.760.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.599.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.346"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.053.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "60.773.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.600.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "355.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.714.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0841"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +9.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  +2.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("E35").Value = "  +4.25%  "
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.918"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.907"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.33%  "
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "296.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Hedera"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0558"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("E49").Value = "  +2.02%  "
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.28%  "
